$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D receive plain numeric-looking text (e.g. "311.11").
# Excel would normally auto-convert such text to a real number, which would
# corrupt the exact original string (trailing zeros, precision, etc.), since
# the source workbook stores these as plain text. Force column D to Text format
# before writing the new values, then clear the temporary formatting so the
# cells end up with no explicit style, matching the original file layout.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '45.464.13'
$ws.Range('E2').Value = '  +6.30%  '

$ws.Range('D3').Value = '2.368.79'
$ws.Range('E3').Value = '  +2.38%  '

$ws.Range('E4').Value = '  +0.17%  '

$ws.Range('D5').Value = '110.55'
$ws.Range('E5').Value = '  +3.19%  '

$ws.Range('D6').Value = '311.11'
$ws.Range('E6').Value = '  -0.30%  '

$ws.Range('D7').Value = '0.630'
$ws.Range('E7').Value = '  +0.43%  '

$ws.Range('E8').Value = '  -0.13%  '

$ws.Range('D9').Value = '0.619'
$ws.Range('E9').Value = '  +1.78%  '

$ws.Range('D10').Value = '41.40'
$ws.Range('E10').Value = '  +3.10%  '

$ws.Range('D11').Value = '0.0921'
$ws.Range('E11').Value = '  +0.52%  '

$ws.Range('D12').Value = '8.51'
$ws.Range('E12').Value = '  +1.68%  '

$ws.Range('E13').Value = '  +1.86%  '

$ws.Range('D14').Value = '0.988'
$ws.Range('E14').Value = '  -0.86%  '

$ws.Range('D15').Value = '2.729.10'
$ws.Range('E15').Value = '  +2.43%  '

$ws.Range('D16').Value = '15.44'
$ws.Range('E16').Value = '  +0.25%  '

$ws.Range('D17').Value = '2.365.68'
$ws.Range('E17').Value = '  +2.34%  '

$ws.Range('D18').Value = '45.393.78'
$ws.Range('E18').Value = '  +6.15%  '

$ws.Range('D19').Value = '7.34'
$ws.Range('E19').Value = '  -1.89%  '

$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').Value = '13.99'
$ws.Range('E20').Value = '  +7.09%  '

$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.0000107'
$ws.Range('E21').Value = '  +0.92%  '

$ws.Range('D22').Value = '73.51'
$ws.Range('E22').Value = '  -0.14%  '

$ws.Range('D23').Value = '3.46'
$ws.Range('E23').Value = '  -0.97%  '

$ws.Range('D24').Value = '259.64'
$ws.Range('E24').Value = '  -2.77%  '

$ws.Range('D25').Value = '2.31'
$ws.Range('E25').Value = '  +2.21%  '

$ws.Range('E26').Value = '  -0.45%  '

$ws.Range('D27').Value = '11.17'
$ws.Range('E27').Value = '  +1.45%  '

$ws.Range('D28').Value = '7.43'
$ws.Range('E28').Value = '  -3.65%  '

$ws.Range('E29').Value = '  +2.55%  '

$ws.Range('D30').Value = '0.0970'
$ws.Range('E30').Value = '  +10.98%  '

$ws.Range('D31').Value = '38.24'
$ws.Range('E31').Value = '  -0.94%  '

$ws.Range('D32').Value = '22.44'
$ws.Range('E32').Value = '  +0.10%  '

$ws.Range('D33').Value = '170.58'

$ws.Range('D34').Value = '2.91'
$ws.Range('E34').Value = '  +6.43%  '

$ws.Range('E35').Value = '  +0.58%  '

$ws.Range('D36').Value = '4.86'
$ws.Range('E36').Value = '  +4.51%  '

$ws.Range('E37').Value = '  +2.46%  '

$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').Value = '3.97'
$ws.Range('E38').Value = '  +8.01%  '

$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').Value = '2.97'
$ws.Range('E39').Value = '  +4.77%  '

$ws.Range('E40').Value = '  +0.30%  '

$ws.Range('D41').Value = '1.75'
$ws.Range('E41').Value = '  +9.57%  '

$ws.Range('D42').Value = '99.74'
$ws.Range('E42').Value = '  -4.16%  '

$ws.Range('D43').Value = '0.234'
$ws.Range('E43').Value = '  +0.42%  '

$ws.Range('D44').Value = '70.16'
$ws.Range('E44').Value = '  -1.63%  '

$ws.Range('D45').Value = '12.89'
$ws.Range('E45').Value = '  +2.65%  '

$ws.Range('E46').Value = '  +0.05%  '

$ws.Range('D47').Value = '83.15'
$ws.Range('E47').Value = '  +8.13%  '

$ws.Range('D48').Value = '113.28'
$ws.Range('E48').Value = '  +0.46%  '

$ws.Range('D49').Value = '9.25'
$ws.Range('E49').Value = '  +4.42%  '

$ws.Range('D50').Value = '5.51'
$ws.Range('E50').Value = '  +4.92%  '

$ws.Range('D51').Value = '1.666.93'
$ws.Range('E51').Value = '  +0.07%  '

$ws.Range('D2:D51').ClearFormats()
